# Updated cryptos list on Sat Jun 15 09:42:43 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $Text) {
    # Prefix with an apostrophe so Excel stores the literal text instead of
    # auto-coercing numeric-looking strings (e.g. "606.69") into numbers,
    # then reset the cell style so no stray quote-prefix formatting lingers.
    $Worksheet.Range($Address).Value = "'" + $Text
    $Worksheet.Range($Address).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "66.153.43"
Set-TextValue $ws "E2" "  -1.28%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "3.524.48"
Set-TextValue $ws "E3" "  +0.15%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  -0.01%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "606.69"
Set-TextValue $ws "E5" "  -0.24%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "143.33"
Set-TextValue $ws "E6" "  -3.16%  "

# Row 7 - LidoStakedEther
Set-TextValue $ws "D7" "3.522.84"
Set-TextValue $ws "E7" "  +0.11%  "

# Row 8 - USDC
Set-TextValue $ws "E8" "  -0.06%  "

# Row 9 - XRP
Set-TextValue $ws "E9" "  +0.17%  "

# Row 10 - was Toncoin, now Dogecoin
Set-TextValue $ws "B10" "Dogecoin"
Set-TextValue $ws "C10" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws "D10" "0.136"
Set-TextValue $ws "E10" "  -4.59%  "

# Row 11 - was Dogecoin, now Toncoin
Set-TextValue $ws "B11" "Toncoin"
Set-TextValue $ws "C11" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D11" "8.05"
Set-TextValue $ws "E11" "  +1.42%  "

# Row 12 - Cardano
Set-TextValue $ws "D12" "0.410"
Set-TextValue $ws "E12" "  -2.89%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D13" "4.120.65"
Set-TextValue $ws "E13" "  +0.13%  "

# Row 14 - ShibaInu
Set-TextValue $ws "D14" "0.0000206"
Set-TextValue $ws "E14" "  -5.16%  "

# Row 15 - Avalanche
Set-TextValue $ws "D15" "30.24"
Set-TextValue $ws "E15" "  -5.39%  "

# Row 16 - WrappedEther
Set-TextValue $ws "D16" "3.523.73"
Set-TextValue $ws "E16" "  +0.22%  "

# Row 17 - WrappedBTC
Set-TextValue $ws "D17" "66.246.02"
Set-TextValue $ws "E17" "  -1.20%  "

# Row 18 - TRON
Set-TextValue $ws "E18" "  -0.67%  "

# Row 19 - Uniswap
Set-TextValue $ws "D19" "10.86"
Set-TextValue $ws "E19" "  +1.38%  "

# Row 20 - Polkadot
Set-TextValue $ws "D20" "6.20"
Set-TextValue $ws "E20" "  -3.42%  "

# Row 21 - Chainlink
Set-TextValue $ws "D21" "14.91"
Set-TextValue $ws "E21" "  -2.71%  "

# Row 22 - BitcoinCash
Set-TextValue $ws "D22" "425.48"
Set-TextValue $ws "E22" "  -2.90%  "

# Row 23 - Polygon
Set-TextValue $ws "D23" "0.600"
Set-TextValue $ws "E23" "  -1.55%  "

# Row 24 - Litecoin
Set-TextValue $ws "D24" "78.66"
Set-TextValue $ws "E24" "  -0.77%  "

# Row 25 - WrappedeETH
Set-TextValue $ws "D25" "3.666.72"
Set-TextValue $ws "E25" "  +0.17%  "

# Row 26 - Dai
Set-TextValue $ws "E26" "  -0.03%  "

# Row 27 - PEPE
Set-TextValue $ws "E27" "  -2.27%  "

# Row 28 - was InternetComputer(DFINITY), now RenderToken
Set-TextValue $ws "B28" "RenderToken"
Set-TextValue $ws "C28" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D28" "7.99"
Set-TextValue $ws "E28" "  -3.93%  "

# Row 29 - was RenderToken, now InternetComputer(DFINITY)
Set-TextValue $ws "B29" "InternetComputer(DFINITY)"
Set-TextValue $ws "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D29" "9.17"
Set-TextValue $ws "E29" "  -6.19%  "

# Row 30 - PancakeSwap
Set-TextValue $ws "D30" "2.46"
Set-TextValue $ws "E30" "  -1.85%  "

# Row 31 - Binance-PegBSC-USD
Set-TextValue $ws "E31" "  +0.37%  "

# Row 32 - Kaspa
Set-TextValue $ws "E32" "  -3.75%  "

# Row 33 - Fetch.AI
Set-TextValue $ws "D33" "1.47"
Set-TextValue $ws "E33" "  -6.60%  "

# Row 34 - EthereumClassic
Set-TextValue $ws "D34" "25.23"
Set-TextValue $ws "E34" "  -0.91%  "

# Row 35 - RenzoRestakedETH
Set-TextValue $ws "D35" "3.515.83"
Set-TextValue $ws "E35" "  +0.10%  "

# Row 36 - USDe
Set-TextValue $ws "E36" "  -0.05%  "

# Row 37 - ImmutableX
Set-TextValue $ws "D37" "1.74"
Set-TextValue $ws "E37" "  -3.51%  "

# Row 38 - Aptos
Set-TextValue $ws "D38" "7.82"
Set-TextValue $ws "E38" "  -2.58%  "

# Row 39 - NEARProtocol
Set-TextValue $ws "D39" "5.59"
Set-TextValue $ws "E39" "  -5.93%  "

# Row 40 - FirstDigitalUSD
Set-TextValue $ws "E40" "  -0.06%  "

# Row 41 - Monero
Set-TextValue $ws "D41" "172.37"
Set-TextValue $ws "E41" "  -0.06%  "

# Row 42 - Hedera
Set-TextValue $ws "D42" "0.0854"
Set-TextValue $ws "E42" "  -4.42%  "

# Row 43 - Filecoin
Set-TextValue $ws "E43" "  -4.85%  "

# Row 44 - Mantle
Set-TextValue $ws "E44" "  -0.59%  "

# Row 45 - Stacks
Set-TextValue $ws "E45" "  -9.32%  "

# Row 46 - OKB
Set-TextValue $ws "D46" "45.27"
Set-TextValue $ws "E46" "  -1.64%  "

# Row 47 - InjectiveProtocol
Set-TextValue $ws "D47" "25.92"
Set-TextValue $ws "E47" "  -7.87%  "

# Row 48 - ONDO
Set-TextValue $ws "D48" "1.20"
Set-TextValue $ws "E48" "  -6.61%  "

# Row 49 - dogwifhat
Set-TextValue $ws "D49" "2.40"
Set-TextValue $ws "E49" "  -2.63%  "

# Row 50 - Cosmos
Set-TextValue $ws "D50" "7.16"
Set-TextValue $ws "E50" "  -4.15%  "

# Row 51 - SuiNetwork
Set-TextValue $ws "D51" "0.944"
Set-TextValue $ws "E51" "  -4.55%  "
